# fix: import excel, generate QR codes, and render PDF
#
# Replaces the placeholder 1..8 demo grid with a real imported dataset
# (id, nominal, nama) so downstream QR/PDF generation has something to
# work from.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old placeholder grid (A1:F8) entirely.
$ws.UsedRange.Clear()

# Row 1: id 1, nominal 141341, nama "harmanto"
$ws.Cells.Item(1, 1).Value = 1
$ws.Cells.Item(1, 2).Value = 141341
$ws.Cells.Item(1, 3).Value = "harmanto"

# Row 2: id 16, nominal 3030000, nama "Khumaidah"
$ws.Cells.Item(2, 1).Value = 16
$ws.Cells.Item(2, 2).Value = 3030000
$ws.Cells.Item(2, 3).Value = "Khumaidah"

# Size the name column to fit the imported text (best-fit, like Excel does
# automatically when a pasted/imported column is auto-sized).
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(3).ColumnWidth = 9.64

# Leave the selection where Excel would land after filling the table.
$ws.Range("C3").Select() | Out-Null

# Restore the (smaller, re-docked) workbook window the file was saved from.
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 10620
    $win.Top = 0
    $win.Width = 9870
    $win.Height = 10920
} catch {
    # Window geometry is cosmetic/host-specific; ignore if unsupported.
}
